$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Map of variable name (column B) -> new description (column L) text.
# These correspond to the shared-string normalization/edits in the diff:
# whitespace/newlines collapsed to single spaces and a few casing tweaks.
$updates = @{
    "icd10_mapped"        = "ICD10 GM diagnosis code mapped A = 1, B = 2, C = 3, D = 4, e.g.: A01.9 = 101.9, C50.1 = 350.1 or D41.9 = 441.9"
    "icd10_grouped"       = "ICD10 GM diagnosis code grouped to parent code, e.g. A01.1 and A01.9 both belong to group 101 (remove decimal from icd10_mapped)"
    "icd10_entity"        = "entities of resulting icd10 groups, see utils"
    "date_diagnosis"      = "date of diagnosis"
    "date_diagnosis_year" = "Year of Diagnosis"
    "age_group_small"     = "age groups mapped as follows: 0 (0-14), 1 (15-19), 2 (20-24), 3 (25-29), 4 (30-34), 5 (35-39), 6 (40-44), 7 (45-49), 8 (50-54), 9 (55-59), 10 (60-64), 11 (65-69), 12 (70-74), 13 (75-79), 14 (80-84), and 15 (85+)"
    "age_group_large"     = "age groups mapped as follows: 0 (0-10), 1 (11-20), 2 (21-30), 3 (31-40), 4 (41-50), 5 (51-60), 6 (61-70), 7 (71-80), 8 (81-90), and 9 (90+)."
    "gender_mapped"       = "Gender mapped: 0 = None, 1 = female, 2 = male, 3 = other/diverse"
    "postal_code"         = "postal code"
}

$rowCount = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $rowCount; $r++) {
    $name = $ws.Cells.Item($r, 2).Value2
    if ($updates.ContainsKey($name)) {
        $ws.Cells.Item($r, 12).Value = $updates[$name]
    }
}
